# This script rewrites numeric "count" cells (column B on the data sheets,
# and A2 on the Overall sheet) as text strings (matching how the source
# report now renders them, e.g. "2,792" instead of 2792), fixes the two
# all-zero County rows to match the new percentage/currency text styling,
# and appends a new "Total" row to the County sheet.
#
# Note: logic is written inline (no helper functions) throughout this
# script, since calling a user-defined function many times in a loop in
# this host triggers intermittent COM marshalling errors.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 numeric 2792 -> text "2,792"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
$cell = $wsOverall.Cells.Item(2, 1)
$cell.NumberFormat = "@"
$cell.Value = "2,792"
$cell.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "County": convert column B (rows 2-132) numeric filer counts to
# plain text of the same number.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
for ($r = 2; $r -le 132; $r++) {
    $cell = $wsCounty.Cells.Item($r, 2)
    $num = $cell.Value()
    $cell.NumberFormat = "@"
    $cell.Value = [string]$num
    $cell.ClearFormats()
}

# Rows 133 (Manassas Park city) and 134 (King George County) were all
# zeros stored inconsistently; restyle them to match the percentage /
# currency text formatting used elsewhere on the sheet.
foreach ($r in 133, 134) {
    $cell = $wsCounty.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = "0.00%"
    $cell.ClearFormats()

    $cell = $wsCounty.Cells.Item($r, 3)
    $cell.NumberFormat = "@"
    $cell.Value = '$0'
    $cell.ClearFormats()

    $cell = $wsCounty.Cells.Item($r, 4)
    $cell.NumberFormat = "@"
    $cell.Value = "0.00%"
    $cell.ClearFormats()

    $cell = $wsCounty.Cells.Item($r, 5)
    $cell.NumberFormat = "@"
    $cell.Value = "0.00%"
    $cell.ClearFormats()

    $cell = $wsCounty.Cells.Item($r, 6)
    $cell.NumberFormat = "@"
    $cell.Value = "0.00%"
    $cell.ClearFormats()
}

# Append new row 135: state-wide "Total" row (mirrors the Total rows
# already present on the other breakdown sheets).
$wsCounty.Cells.Item(135, 1).Value = "Total"

$cell = $wsCounty.Cells.Item(135, 2)
$cell.NumberFormat = "@"
$cell.Value = "2,792"
$cell.ClearFormats()

$cell = $wsCounty.Cells.Item(135, 3)
$cell.NumberFormat = "@"
$cell.Value = '$4,571,835,868'
$cell.ClearFormats()

$cell = $wsCounty.Cells.Item(135, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.66%"
$cell.ClearFormats()

$cell = $wsCounty.Cells.Item(135, 5)
$cell.NumberFormat = "@"
$cell.Value = "-9.96%"
$cell.ClearFormats()

$cell = $wsCounty.Cells.Item(135, 6)
$cell.NumberFormat = "@"
$cell.Value = "64.29%"
$cell.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "Congressional District": convert column B (rows 2-13) numeric
# values to text. Row 13 is the "Total" row (2792 -> "2,792"); the rest
# keep their plain digits.
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 12; $r++) {
    $cell = $wsCd.Cells.Item($r, 2)
    $num = $cell.Value()
    $cell.NumberFormat = "@"
    $cell.Value = [string]$num
    $cell.ClearFormats()
}
$cell = $wsCd.Cells.Item(13, 2)
$cell.NumberFormat = "@"
$cell.Value = "2,792"
$cell.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "Size": convert column B (rows 2-8) numeric values to text. Row 8
# is the "Total" row (2792 -> "2,792").
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 7; $r++) {
    $cell = $wsSize.Cells.Item($r, 2)
    $num = $cell.Value()
    $cell.NumberFormat = "@"
    $cell.Value = [string]$num
    $cell.ClearFormats()
}
$cell = $wsSize.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = "2,792"
$cell.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "Subsector": convert column B (rows 2-13) numeric values to text.
# Row 13 is the "Total" row (2792 -> "2,792").
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 12; $r++) {
    $cell = $wsSub.Cells.Item($r, 2)
    $num = $cell.Value()
    $cell.NumberFormat = "@"
    $cell.Value = [string]$num
    $cell.ClearFormats()
}
$cell = $wsSub.Cells.Item(13, 2)
$cell.NumberFormat = "@"
$cell.Value = "2,792"
$cell.ClearFormats()
